$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "pop_sq_mile_5mi"

$ws.Range("R2").Value = 41.6216216216216
$ws.Range("S2").Value = 0.475675675675676

$ws.Range("R3").Value = 46.6666666666667
$ws.Range("S3").Value = 0.5

$ws.Range("R4").Value = 32.2857142857143
$ws.Range("S4").Value = 0.3

$ws.Range("R5").Value = 30.9625668449198
$ws.Range("S5").Value = 0.383422459893048

$ws.Range("R6").Value = 29.4545454545455
$ws.Range("S6").Value = 0.377272727272727

$ws.Range("R7").Value = 30
$ws.Range("S7").Value = 0.336363636363636

$ws.Range("R8").Value = 61.4285714285714
$ws.Range("S8").Value = 0.45

$ws.Range("R9").Value = 61.8181818181818
$ws.Range("S9").Value = 0.463636363636364

$ws.Range("R10").Value = 20
$ws.Range("S10").Value = 0.2

$ws.Range("R11").Value = 19.047619047619
$ws.Range("S11").Value = 0.185714285714286
